$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D header: "WebExcel" -> "ExDataExcel"
$ws.Range("D1").Value = "ExDataExcel"

# Column D data rows: "...WebData.xlsx" -> "...ExcelData.xlsx"
$ws.Range("D2").Value = "TC03_CDSValidation_by_ParticipantID - 7_ExcelData.xlsx"
$ws.Range("D3").Value = "TC03_CDSValidation_by_ParticipantID - 7_ExcelData.xlsx"
$ws.Range("D4").Value = "TC03_CDSValidation_by_ParticipantID - 7_ExcelData.xlsx"

# Update the active cell selection to D4 (matches final saved view state)
$ws.Range("D4").Select() | Out-Null
